$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 (2020-09-14 Mon) ---
# Start/end time worked, which derives D25 (=C25-B25) and F25 (=SUM(D21:D25)) automatically.
$ws.Range("B25").Value = 10
$ws.Range("C25").Value = 14
$ws.Range("E25").Value = "Bottom NavBar styling, building realmlist to ListView with unique elements from array"

# --- Row 26 (2020-09-17 Thu) ---
# Updated start/end time, which derives D26 (=C26-B26) and feeds F30 (=SUM(D26:D30)) automatically.
$ws.Range("B26").Value = 9
$ws.Range("C26").Value = 16
$ws.Range("E26").Value = "Input fields for character and realm, realm selection, fetching the values of the input fields, some modifications for RaiderIO API"

# Update the view's active selection to reflect where the author ended up editing.
$ws.Range("E28").Select()
